# Revert "Merging 0.1.8 w VitalSigns"
#
# - Rename sheet "Include #0" -> "Include from SNOMED CT"
# - Metadata sheet: roll several property values back to their pre-merge
#   values, and delete the "Jurisdiction" row (shifting rows below it up).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInc  = $wb.Worksheets.Item("Include #0")

# Rename the Include sheet.
$wsInc.Name = "Include from SNOMED CT"

# Roll back simple property values on the Metadata sheet.
$wsMeta.Range("B3").Value = "0.1.6"
$wsMeta.Range("B6").Value = "active"
$wsMeta.Range("B8").Value = "2023-05-05T10:50:04-05:00"
$wsMeta.Range("B10").Value = "No display for ContactDetail"
$wsMeta.Range("B11").Value = "No display for ContactDetail"

# Delete the entire "Jurisdiction" row (row 12); everything below shifts up.
$wsMeta.Range("A12").EntireRow.Delete()
